$wb = $excel.ActiveWorkbook

# --- concepts sheet: add a new "Thing" concept row ---
$concepts = $wb.Worksheets.Item("concepts")
$concepts.Activate()
$concepts.Range("A2").Value = "default"
$concepts.Range("B2").Value = "Thing"
$concepts.Range("C2").Value = "Thing summary"
$concepts.Range("D2").Value = "Thing description."
$concepts.Range("D2").Select() | Out-Null

# --- elements sheet: add a new "hasThing" element row ---
$elements = $wb.Worksheets.Item("elements")
$elements.Activate()
$elements.Range("A2").Value = "default"
$elements.Range("B2").Value = "hasThing"
$elements.Range("C2").Value = "hasThing summary"
$elements.Range("B5").Select() | Out-Null

# --- packages sheet: widen column B to fit the "Default package" text ---
$packages = $wb.Worksheets.Item("packages")
$packages.Columns.Item(2).ColumnWidth = 16.2
